# Team_Results.xlsx update — "Ran through latest numbers / Cleaned up the folder as well."
#
# This applies the season-stats refresh captured in the diff:
#  - updated "Last POS" delta column (Y) for several clubs (new rows gained
#    a value, a couple of rows lost theirs, most changed)
#  - the per-column "best position" marker row (row 24) picked up new
#    winners on a few columns, including two brand-new columns (O, X)
#  - the highlight color used on that marker row moved from yellow to a
#    dark red, and one cell (D24) lost its highlight entirely
#  - the view was scrolled back to the left edge and the selection left on H25
#  - the header row got a little shorter and the page setup got a paper size

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Last POS" (column Y) numbers refreshed for the new results ---
$ws.Range("Y2").ClearContents()
$ws.Range("Y3").Value = 2.75
$ws.Range("Y6").Value = 2.59
$ws.Range("Y7").Value = 2.57
$ws.Range("Y8").Value = 3.52
$ws.Range("Y9").Value = 2.46
$ws.Range("Y10").ClearContents()
$ws.Range("Y11").Value = 2.99
$ws.Range("Y18").Value = 2.24
$ws.Range("Y19").Value = 2.89
$ws.Range("Y22").Value = 2.26

# --- row 24 "best position per column" labels ---
$ws.Range("K24").Value = "LW"
$ws.Range("N24").Value = "CD"
$ws.Range("O24").Value = "IM"
$ws.Range("R24").Value = "LB"
$ws.Range("U24").Value = "CD"
$ws.Range("X24").Value = "IM"

# --- highlight formatting on row 24: yellow -> dark red, D24 cleared ---
# (dark red FFC00000 == R192 G0 B0; Excel OLE colors are 0x00BBGGRR)
$darkRed = 192
$ws.Range("D24").ClearFormats()
$ws.Range("C24").Interior.Color = $darkRed
$ws.Range("E24").Interior.Color = $darkRed
$ws.Range("I24").Interior.Color = $darkRed

# --- header row height trimmed a bit ---
$ws.Range("A1:Y1").RowHeight = 75

# --- page setup: portrait, paper size 138 ---
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 138

# --- view: scrolled back to the left edge, selection left on H25 ---
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("H25").Select() | Out-Null
